$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 46, pushing existing rows 46-55 down to 47-56.
$ws.Rows.Item(46).Insert()

# Populate the new row 46 with the values from the diff. Columns A,B,C,E,F,G,H,I,J
# carry the same constant values as the surrounding rows in this block.
$ws.Cells.Item(46, 1).Value = 10
$ws.Cells.Item(46, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(46, 3).Value = "La Araucanía"
$ws.Cells.Item(46, 4).Value = 45027
$ws.Cells.Item(46, 4).NumberFormat = $ws.Cells.Item(47, 4).NumberFormat
$ws.Cells.Item(46, 5).Value = 9
$ws.Cells.Item(46, 6).Value = "Fruta"
$ws.Cells.Item(46, 7).Value = 100107
$ws.Cells.Item(46, 8).Value = "Otros"
$ws.Cells.Item(46, 9).Value = 100107001
$ws.Cells.Item(46, 10).Value = "Caqui"
$ws.Cells.Item(46, 11).Value = "Fuyu"
$ws.Cells.Item(46, 12).Value = "Primera"
$ws.Cells.Item(46, 13).Value = 55
$ws.Cells.Item(46, 14).Value = 24000
$ws.Cells.Item(46, 15).Value = 24000
$ws.Cells.Item(46, 16).Value = 24000
$ws.Cells.Item(46, 17).Value = "`$/bandeja 15 kilos granel"
$ws.Cells.Item(46, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(46, 19).Value = 1600
$ws.Cells.Item(46, 20).Value = 15
